$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update moment analysis data (PH columns I-L) for rows 12, 13, 14, 16
$ws.Range("I12").Value = 0.1741306222569814
$ws.Range("J12").Value = 0.007661405384502213
$ws.Range("K12").Value = -0.1494289167025337
$ws.Range("L12").Value = 1.834056543320003

$ws.Range("I13").Value = 0.1732893829452975
$ws.Range("J13").Value = 0.006388042904982789
$ws.Range("K13").Value = -0.1320200159195386
$ws.Range("L13").Value = 2.192210077531122

$ws.Range("I14").Value = 0.1897922458443881
$ws.Range("J14").Value = 0.00710953225032518
$ws.Range("K14").Value = -0.3236641155507869
$ws.Range("L14").Value = 2.042420822161466

$ws.Range("I16").Value = 0.2410106334485396
$ws.Range("J16").Value = 0.02458410332041821
$ws.Range("K16").Value = 1.171992255373215
$ws.Range("L16").Value = 3.796393522117083
